$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 323 (pushes existing rows 323-359 down to 324-360),
#    inheriting formatting from the row above (row 322).
$ws.Rows("323:323").Insert()

# 2) Populate the newly inserted row with the new GenTable entry.
$ws.Range("A323").Value = "L9-報表作業"
$ws.Range("B323").Value = "YearlyHouseLoanIntCheck"
$ws.Range("C323").Value = "每年房屋擔保借款繳息檢核檔"
$ws.Range("D323").Formula = '=HYPERLINK("[\\192.168.10.16\St1Share(NAS)\SKL\DB\GenTables\L9-報表作業\YearlyHouseLoanIntCheck.xlsx]DBD!A1", "連結")'
$ws.Range("E323").Value = "2022年03月23日 17:47:49"

# 3) Update the three "last modified" timestamps that changed elsewhere in the table.
$ws.Range("E37").Value = "2022年03月23日 18:37:00"
$ws.Range("E68").Value = "2022年03月23日 17:25:59"
$ws.Range("E85").Value = "2022年03月24日 09:32:37"
